$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 17 data
$ws.Range("A17").Value = "PT Bima (Sunco)"
$ws.Range("B17").Value = 45515
$ws.Range("C17").Value = 46610
$ws.Range("D17").Value = ""
$ws.Range("E17").Value = 3
$ws.Range("F17").Value = 235000000
$ws.Range("G17").Value = 235000000
$ws.Range("H17").Value = "Full Lease Upfront"
$ws.Range("I17").Value = ""
$ws.Range("J17").Value = ""

# Copy style (number format) from row 16's date cells to row 17's date cells
$ws.Range("B16").Copy() | Out-Null
$ws.Range("B17").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$ws.Range("C16").Copy() | Out-Null
$ws.Range("C17").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
